$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("多氟多", "闻泰科技", "闻泰科技")
    3  = @("海马汽车", "多氟多", "合富中国")
    4  = @("天际股份", "合富中国", "平潭发展")
    5  = @("闻泰科技", "特变电工", "多氟多")
    6  = @("合富中国", "平潭发展", "兰石重装")
    7  = @("澄星股份", "天赐材料", "天赐材料")
    8  = @("平潭发展", "海马汽车", "雪人集团")
    9  = @("方正电机", "天际股份", "方正电机")
    10 = @("特变电工", "方正电机", "海马汽车")
    11 = @("天赐材料", "兰石重装", "特变电工")
    12 = @("兰石重装", "海陆重工", "摩恩电气")
    13 = @("合盛硅业", "永太科技", "隆基绿能")
    14 = @("海陆重工", "澄星股份", "澄星股份")
    15 = @("摩恩电气", "东方财富", "三花智控")
    16 = @("中国西电", "江苏国泰", "海天股份")
    17 = @("深圳新星", "东岳硅材", "中国西电")
    18 = @("隆基绿能", "中国西电", "盈新发展")
    19 = @("雪人集团", "贵州茅台", "漳州发展")
    20 = @("永太科技", "合盛硅业", "中毅达")
    21 = @("洲际油气", "摩恩电气", "粤传媒")
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("A$row").Value = $values[0]
    $ws.Range("B$row").Value = $values[1]
    $ws.Range("C$row").Value = $values[2]
}
